$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "check_is_peak_area" column (I)
# so it becomes check_ion_ratio | check_ion_ratio_std | check_is_peak_area | check_spike
$ws.Columns("I:I").Insert()

# New column header
$ws.Range("I1").Value = "check_ion_ratio_std"

# New column values for the first 14 compounds (the "6,7,10,11,a,b,c,d" / std-curve group)
$ws.Range("I2:I15").Value = "1,2,3,4,5"

# The glucuronide compounds (rows 16-21) leave the new column blank
$ws.Range("I16:I21").ClearContents()

# Existing data updates that accompany the new column
$ws.Range("D2:D15").Value = "6,7,10,11,a,b,c,d"
$ws.Range("E16:E21").Value = "8,10,11,a,c"

# Restore the selection to match the author's saved cursor position
$ws.Range("F24").Select()
